$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDFVALUES")

# Row 2 - update HtmlPath, ScreenShotsPath, TesterName values for new (negative) test run
$ws.Range("A2").Value = "C:/Users/Dell/Downloads/pAInITe-master/pAInITe-master/test-output/Default%20suite/Default%20test.html"
$ws.Range("B2").Value = "C:\Users\Dell\Downloads\pAInITe-master\pAInITe-master\TestcaseScreenshots\"
$ws.Range("C2").Value = "Deepika"
$ws.Range("E2").Value = "C:\Users\Dell\Downloads\pAInITe-master\pAInITe-master\Logo\M10logo.png"

# B2 now wraps and centers vertically
$ws.Range("B2").WrapText = $true
$ws.Range("B2").VerticalAlignment = -4108

# Row grows to fit the now-wrapped text
$ws.Rows.Item(2).RowHeight = 60

# Clear out the old screenshots-path row entirely (row becomes empty -> vanishes)
$ws.Range("B6").ClearContents()

# Clear the leftover sample row's data, keeping its formatting
$ws.Range("A8:E8").ClearContents()
$ws.Rows.Item(8).AutoFit()

$ws.Range("A6:E8").Select()
